# Update "想去人数" (want-to-go count) / "最低票价" (min price) figures
# across all four sheets to match the refreshed scrape output
# (gh-pages rebuild at 456a3b4).
$wb = $excel.ActiveWorkbook


# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
if ($ws.Range("F3").Value() -ne 3467) { throw "Unexpected value in 展览!F3: expected 3467" }
$ws.Range("F3").Value = 3468
if ($ws.Range("F4").Value() -ne 368) { throw "Unexpected value in 展览!F4: expected 368" }
$ws.Range("F4").Value = 369
if ($ws.Range("F5").Value() -ne 8125) { throw "Unexpected value in 展览!F5: expected 8125" }
$ws.Range("F5").Value = 8126
if ($ws.Range("F8").Value() -ne 2105) { throw "Unexpected value in 展览!F8: expected 2105" }
$ws.Range("F8").Value = 2106
if ($ws.Range("G10").Value() -ne 58) { throw "Unexpected value in 展览!G10: expected 58" }
$ws.Range("G10").Value = 68
if ($ws.Range("F12").Value() -ne 506) { throw "Unexpected value in 展览!F12: expected 506" }
$ws.Range("F12").Value = 507
if ($ws.Range("F16").Value() -ne 1150) { throw "Unexpected value in 展览!F16: expected 1150" }
$ws.Range("F16").Value = 1151
if ($ws.Range("F18").Value() -ne 721) { throw "Unexpected value in 展览!F18: expected 721" }
$ws.Range("F18").Value = 722
if ($ws.Range("F22").Value() -ne 416) { throw "Unexpected value in 展览!F22: expected 416" }
$ws.Range("F22").Value = 418
if ($ws.Range("F24").Value() -ne 4723) { throw "Unexpected value in 展览!F24: expected 4723" }
$ws.Range("F24").Value = 4879
if ($ws.Range("F26").Value() -ne 50689) { throw "Unexpected value in 展览!F26: expected 50689" }
$ws.Range("F26").Value = 50991
if ($ws.Range("F27").Value() -ne 4032) { throw "Unexpected value in 展览!F27: expected 4032" }
$ws.Range("F27").Value = 4041
if ($ws.Range("F29").Value() -ne 991) { throw "Unexpected value in 展览!F29: expected 991" }
$ws.Range("F29").Value = 992
if ($ws.Range("F30").Value() -ne 759) { throw "Unexpected value in 展览!F30: expected 759" }
$ws.Range("F30").Value = 761
if ($ws.Range("F31").Value() -ne 337) { throw "Unexpected value in 展览!F31: expected 337" }
$ws.Range("F31").Value = 343
if ($ws.Range("F33").Value() -ne 834) { throw "Unexpected value in 展览!F33: expected 834" }
$ws.Range("F33").Value = 835
if ($ws.Range("F35").Value() -ne 567) { throw "Unexpected value in 展览!F35: expected 567" }
$ws.Range("F35").Value = 568
if ($ws.Range("F38").Value() -ne 1) { throw "Unexpected value in 展览!F38: expected 1" }
$ws.Range("F38").Value = 2
if ($ws.Range("F39").Value() -ne 827) { throw "Unexpected value in 展览!F39: expected 827" }
$ws.Range("F39").Value = 829
if ($ws.Range("F40").Value() -ne 1019) { throw "Unexpected value in 展览!F40: expected 1019" }
$ws.Range("F40").Value = 1023
if ($ws.Range("F42").Value() -ne 157) { throw "Unexpected value in 展览!F42: expected 157" }
$ws.Range("F42").Value = 159
if ($ws.Range("F47").Value() -ne 91) { throw "Unexpected value in 展览!F47: expected 91" }
$ws.Range("F47").Value = 92

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
if ($ws.Range("F11").Value() -ne 110) { throw "Unexpected value in 演出!F11: expected 110" }
$ws.Range("F11").Value = 111
if ($ws.Range("F14").Value() -ne 34) { throw "Unexpected value in 演出!F14: expected 34" }
$ws.Range("F14").Value = 36
if ($ws.Range("F15").Value() -ne 84) { throw "Unexpected value in 演出!F15: expected 84" }
$ws.Range("F15").Value = 85
if ($ws.Range("F17").Value() -ne 27) { throw "Unexpected value in 演出!F17: expected 27" }
$ws.Range("F17").Value = 28
if ($ws.Range("F18").Value() -ne 154) { throw "Unexpected value in 演出!F18: expected 154" }
$ws.Range("F18").Value = 155
if ($ws.Range("F19").Value() -ne 7302) { throw "Unexpected value in 演出!F19: expected 7302" }
$ws.Range("F19").Value = 7304
if ($ws.Range("F20").Value() -ne 67) { throw "Unexpected value in 演出!F20: expected 67" }
$ws.Range("F20").Value = 69
if ($ws.Range("F28").Value() -ne 104) { throw "Unexpected value in 演出!F28: expected 104" }
$ws.Range("F28").Value = 105

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
if ($ws.Range("F4").Value() -ne 2203) { throw "Unexpected value in 本地生活!F4: expected 2203" }
$ws.Range("F4").Value = 2204
if ($ws.Range("F9").Value() -ne 9246) { throw "Unexpected value in 本地生活!F9: expected 9246" }
$ws.Range("F9").Value = 9247
if ($ws.Range("F10").Value() -ne 1521) { throw "Unexpected value in 本地生活!F10: expected 1521" }
$ws.Range("F10").Value = 1524
if ($ws.Range("F11").Value() -ne 145) { throw "Unexpected value in 本地生活!F11: expected 145" }
$ws.Range("F11").Value = 146
if ($ws.Range("F13").Value() -ne 1) { throw "Unexpected value in 本地生活!F13: expected 1" }
$ws.Range("F13").Value = 2
if ($ws.Range("F14").Value() -ne 59) { throw "Unexpected value in 本地生活!F14: expected 59" }
$ws.Range("F14").Value = 62

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
if ($ws.Range("F2").Value() -ne 3467) { throw "Unexpected value in 全部类型!F2: expected 3467" }
$ws.Range("F2").Value = 3468
if ($ws.Range("F3").Value() -ne 2203) { throw "Unexpected value in 全部类型!F3: expected 2203" }
$ws.Range("F3").Value = 2204
if ($ws.Range("F6").Value() -ne 1521) { throw "Unexpected value in 全部类型!F6: expected 1521" }
$ws.Range("F6").Value = 1524
if ($ws.Range("F7").Value() -ne 145) { throw "Unexpected value in 全部类型!F7: expected 145" }
$ws.Range("F7").Value = 146
if ($ws.Range("F10").Value() -ne 2105) { throw "Unexpected value in 全部类型!F10: expected 2105" }
$ws.Range("F10").Value = 2106
if ($ws.Range("F13").Value() -ne 506) { throw "Unexpected value in 全部类型!F13: expected 506" }
$ws.Range("F13").Value = 507
if ($ws.Range("F19").Value() -ne 1150) { throw "Unexpected value in 全部类型!F19: expected 1150" }
$ws.Range("F19").Value = 1151
if ($ws.Range("F20").Value() -ne 721) { throw "Unexpected value in 全部类型!F20: expected 721" }
$ws.Range("F20").Value = 722
if ($ws.Range("F23").Value() -ne 416) { throw "Unexpected value in 全部类型!F23: expected 416" }
$ws.Range("F23").Value = 418
if ($ws.Range("F25").Value() -ne 110) { throw "Unexpected value in 全部类型!F25: expected 110" }
$ws.Range("F25").Value = 111
if ($ws.Range("F27").Value() -ne 991) { throw "Unexpected value in 全部类型!F27: expected 991" }
$ws.Range("F27").Value = 992
if ($ws.Range("F28").Value() -ne 337) { throw "Unexpected value in 全部类型!F28: expected 337" }
$ws.Range("F28").Value = 343
if ($ws.Range("F30").Value() -ne 567) { throw "Unexpected value in 全部类型!F30: expected 567" }
$ws.Range("F30").Value = 568
if ($ws.Range("F31").Value() -ne 84) { throw "Unexpected value in 全部类型!F31: expected 84" }
$ws.Range("F31").Value = 85
if ($ws.Range("F34").Value() -ne 27) { throw "Unexpected value in 全部类型!F34: expected 27" }
$ws.Range("F34").Value = 28
if ($ws.Range("F35").Value() -ne 7302) { throw "Unexpected value in 全部类型!F35: expected 7302" }
$ws.Range("F35").Value = 7304
if ($ws.Range("F36").Value() -ne 827) { throw "Unexpected value in 全部类型!F36: expected 827" }
$ws.Range("F36").Value = 829
if ($ws.Range("F37").Value() -ne 67) { throw "Unexpected value in 全部类型!F37: expected 67" }
$ws.Range("F37").Value = 69
if ($ws.Range("F44").Value() -ne 91) { throw "Unexpected value in 全部类型!F44: expected 91" }
$ws.Range("F44").Value = 92
